{"js": "// Replace the 25 multiplication-problem strings in the single 20x5 table.\n// The table has 5 \"data\" rows (0, 4, 9, 14, 19), each with 5 populated\n// cells (columns 0-4); the remaining rows are blank work-space rows.\n// Because some source strings repeat (e.g. \"127\u00d74=\" appears twice) the\n// mapping below is applied by explicit (row, col) position rather than by\n// a global find/replace, so each occurrence maps to the correct target.\n\nconst replacements = [\n  // [row, col, oldText, newText]\n  [0, 0, \"137\u00d76=\", \"248\u00d72=\"],\n  [0, 1, \"470\u00d74=\", \"644\u00d76=\"],\n  [0, 2, \"362\u00d74=\", \"876\u00d72=\"],\n  [0, 3, \"256\u00d73=\", \"637\u00d75=\"],\n  [0, 4, \"472\u00d75=\", \"883\u00d74=\"],\n\n  [4, 0, \"419\u00d78=\", \"719\u00d75=\"],\n  [4, 1, \"940\u00d72=\", \"451\u00d74=\"],\n  [4, 2, \"618\u00d77=\", \"340\u00d79=\"],\n  [4, 3, \"758\u00d73=\", \"308\u00d76=\"],\n  [4, 4, \"918\u00d73=\", \"123\u00d73=\"],\n\n  [9, 0, \"127\u00d74=\", \"954\u00d78=\"],\n  [9, 1, \"514\u00d76=\", \"929\u00d77=\"],\n  [9, 2, \"259\u00d75=\", \"440\u00d78=\"],\n  [9, 3, \"138\u00d78=\", \"568\u00d73=\"],\n  [9, 4, \"512\u00d78=\", \"602\u00d74=\"],\n\n  [14, 0, \"878\u00d78=\", \"554\u00d73=\"],\n  [14, 1, \"272\u00d78=\", \"308\u00d77=\"],\n  [14, 2, \"470\u00d75=\", \"248\u00d78=\"],\n  [14, 3, \"127\u00d74=\", \"360\u00d72=\"],\n  [14, 4, \"614\u00d79=\", \"370\u00d75=\"],\n\n  [19, 0, \"848\u00d78=\", \"919\u00d79=\"],\n  [19, 1, \"185\u00d76=\", \"271\u00d78=\"],\n  [19, 2, \"606\u00d78=\", \"756\u00d73=\"],\n  [19, 3, \"755\u00d78=\", \"771\u00d76=\"],\n  [19, 4, \"927\u00d76=\", \"204\u00d77=\"],\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfor (const [row, col, oldText, newText] of replacements) {\n  const cell = table.getCell(row, col);\n  const cellBody = cell.body;\n  const results = cellBody.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  } else {\n    // Fallback: replace whole cell text if the exact string wasn't found\n    // (keeps the script resilient to minor formatting differences).\n    cellBody.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 multiplication-problem strings in the single 20x5 table.\n# The table has 5 \"data\" rows (1, 5, 10, 15, 20 in 1-based COM indexing),\n# each with 5 populated cells (columns 1-5); the remaining rows are blank\n# work-space rows.\n#\n# Because some source strings repeat (e.g. \"127\u00d74=\" appears twice) we\n# address each target by its exact (row, col) table position rather than\n# a document-wide Find/Replace, and we write only the literal text span\n# (re-deriving it from the cell's own Range) so the existing run/paragraph\n# formatting (font, size, alignment) is left untouched.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nfunction Set-CellText($row, $col, $newText) {\n    $cell = $tbl.Cell($row, $col)\n    $rng = $cell.Range\n    $full = $rng.Text\n    # A table-cell Range.Text includes the trailing cell-mark control\n    # character(s); trim at the first CR (chr 13) to get the real content.\n    $idx = $full.IndexOf([char]13)\n    if ($idx -ge 0) {\n        $len = $idx\n    } else {\n        $len = $full.Length\n    }\n    $textRng = $d.Range($rng.Start, $rng.Start + $len)\n    $textRng.Text = $newText\n}\n\n# Row 1\nSet-CellText 1 1 \"248\u00d72=\"\nSet-CellText 1 2 \"644\u00d76=\"\nSet-CellText 1 3 \"876\u00d72=\"\nSet-CellText 1 4 \"637\u00d75=\"\nSet-CellText 1 5 \"883\u00d74=\"\n\n# Row 5\nSet-CellText 5 1 \"719\u00d75=\"\nSet-CellText 5 2 \"451\u00d74=\"\nSet-CellText 5 3 \"340\u00d79=\"\nSet-CellText 5 4 \"308\u00d76=\"\nSet-CellText 5 5 \"123\u00d73=\"\n\n# Row 10\nSet-CellText 10 1 \"954\u00d78=\"\nSet-CellText 10 2 \"929\u00d77=\"\nSet-CellText 10 3 \"440\u00d78=\"\nSet-CellText 10 4 \"568\u00d73=\"\nSet-CellText 10 5 \"602\u00d74=\"\n\n# Row 15\nSet-CellText 15 1 \"554\u00d73=\"\nSet-CellText 15 2 \"308\u00d77=\"\nSet-CellText 15 3 \"248\u00d78=\"\nSet-CellText 15 4 \"360\u00d72=\"\nSet-CellText 15 5 \"370\u00d75=\"\n\n# Row 20\nSet-CellText 20 1 \"919\u00d79=\"\nSet-CellText 20 2 \"271\u00d78=\"\nSet-CellText 20 3 \"756\u00d73=\"\nSet-CellText 20 4 \"771\u00d76=\"\nSet-CellText 20 5 \"204\u00d77=\"\n"}
